$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Julius Randle"
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "Minnesota Timberwolves"

$ws.Range("A9").Value = "Jonathan Isaac"
$ws.Range("B9").Value = "SF,PF"
$ws.Range("C9").Value = "Orlando Magic"

$ws.Range("A15").Value = "Isaiah Hartenstein"
$ws.Range("B15").Value = "C"
$ws.Range("C15").Value = "Oklahoma City Thunder"

$ws.Range("A16").Value = "Naz Reid"
$ws.Range("B16").Value = "PF,C"
$ws.Range("C16").Value = "Minnesota Timberwolves"

$ws.Range("A17").Value = "Keyonte George"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("C17").Value = "Utah Jazz"

$ws.Range("A18").Value = "Collin Sexton"
$ws.Range("B18").Value = "PG,SG"
$ws.Range("C18").Value = "Utah Jazz"
